$d = $word.ActiveDocument

# Locate the paragraph that holds the trailing copyright/footer notice
# ("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
# pages. Original theme under Creative Commons Attribution").
$copyrightIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Powered by Jekyll*") {
        $copyrightIndex = $i
    }
}

if ($copyrightIndex -gt 2) {
    # The footer notice is preceded by two structural (text-less)
    # paragraphs: a blank "Normal" paragraph and a page-break paragraph.
    # Remove all three together (including their paragraph marks) so the
    # bibliography's last entry is followed directly by the original
    # trailing blank / page-break paragraphs that close the document.
    $startPara = $d.Paragraphs.Item($copyrightIndex - 2)
    $endPara = $d.Paragraphs.Item($copyrightIndex)

    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
